$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before current row 11 ("Are we ready?") to make room
# for the new problems 9, 10, 11. This pushes existing rows 11-12 down to 13-14.
$ws.Range("A11:A12").Insert()

# Fill in the new problem rows (10, 11, 12)
$ws.Range("A10").Value = "9. Since I am the first speaker may be need a bit introduction. So I supposed to add an introduction slide"
$ws.Range("A11").Value = "10. At the introduction do I need to tell the overview of presentation including everyone topics"
$ws.Range("A12").Value = "11. I feel like need to change the order of 1st topic and 2nd topic. Since need to introduce the twitter in the 1st slide Anuradha can get"

# Restore selection as in the edited workbook
$ws.Range("A17").Select()
